$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# B18: ParentID value changes from "17" to "11" (text-formatted cell)
$ws.Range("B18").Value = "11"

# B26: the ParentID cell is removed entirely (was "1")
$ws.Range("B26").Clear()
